# Update cryptos list with latest pulled values (GitHub Actions data refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Some new price figures are plain decimals (e.g. 491.18) that Excel would
# otherwise auto-convert to numbers; the source cells are text, so mark those
# specific cells as Text before writing them to preserve the string type.
$textCells = $excel.Union($ws.Range("D5"), $ws.Range("D6"), $ws.Range("D8"), $ws.Range("D9"), $ws.Range("D11"), $ws.Range("D12"), $ws.Range("D13"), $ws.Range("D16"), $ws.Range("D18"), $ws.Range("D21"), $ws.Range("D23"), $ws.Range("D24"), $ws.Range("D25"), $ws.Range("D26"), $ws.Range("D27"), $ws.Range("D28"), $ws.Range("D30"), $ws.Range("D31"), $ws.Range("D34"), $ws.Range("D36"), $ws.Range("D37"), $ws.Range("D38"), $ws.Range("D40"), $ws.Range("D43"), $ws.Range("D45"), $ws.Range("D49"), $ws.Range("D50"), $ws.Range("D51"))
$textCells.NumberFormat = "@"

$ws.Range("D2").Value = '69.476.82'
$ws.Range("E2").Value = '  +1.86%  '
$ws.Range("D3").Value = '3.947.80'
$ws.Range("E3").Value = '  +0.55%  '
$ws.Range("E4").Value = '  -0.06%  '
$ws.Range("D5").Value = '491.18'
$ws.Range("E5").Value = '  +0.90%  '
$ws.Range("D6").Value = '147.15'
$ws.Range("E6").Value = '  +0.39%  '
$ws.Range("D8").Value = '0.998'
$ws.Range("E8").Value = '  -0.03%  '
$ws.Range("D9").Value = '0.738'
$ws.Range("E9").Value = '  +0.99%  '
$ws.Range("E10").Value = '  +4.13%  '
$ws.Range("D11").Value = '0.0000348'
$ws.Range("E11").Value = '  -3.57%  '
$ws.Range("D12").Value = '43.14'
$ws.Range("E12").Value = '  +0.73%  '
$ws.Range("D13").Value = '10.50'
$ws.Range("E13").Value = '  -1.55%  '
$ws.Range("D14").Value = '4.571.19'
$ws.Range("E14").Value = '  +0.46%  '
$ws.Range("D15").Value = '3.949.74'
$ws.Range("E15").Value = '  +1.05%  '
$ws.Range("D16").Value = '14.26'
$ws.Range("E16").Value = '  -4.53%  '
$ws.Range("E17").Value = '  -0.69%  '
$ws.Range("D18").Value = '19.91'
$ws.Range("E18").Value = '  -1.25%  '
$ws.Range("E19").Value = '  +2.44%  '
$ws.Range("D20").Value = '69.527.52'
$ws.Range("E20").Value = '  +1.78%  '
$ws.Range("D21").Value = '439.63'
$ws.Range("E21").Value = '  -1.81%  '
$ws.Range("E22").Value = '  +2.03%  '
$ws.Range("D23").Value = '14.58'
$ws.Range("E23").Value = '  -1.37%  '
$ws.Range("D24").Value = '89.53'
$ws.Range("E24").Value = '  +0.87%  '
$ws.Range("D25").Value = '12.13'
$ws.Range("E25").Value = '  +12.06%  '
$ws.Range("D26").Value = '3.74'
$ws.Range("E26").Value = '  +3.28%  '
$ws.Range("D27").Value = '11.18'
$ws.Range("E27").Value = '  -2.88%  '
$ws.Range("D28").Value = '37.41'
$ws.Range("E28").Value = '  -4.15%  '
$ws.Range("E29").Value = '  -3.70%  '
$ws.Range("D30").Value = '707.74'
$ws.Range("E30").Value = '  +2.48%  '
$ws.Range("D31").Value = '13.53'
$ws.Range("E31").Value = '  +0.09%  '
$ws.Range("E32").Value = '  +0.16%  '
$ws.Range("E33").Value = '  +0.90%  '
$ws.Range("D34").Value = '0.470'
$ws.Range("E34").Value = '  +25.69%  '
$ws.Range("D35").Value = '0.0₃0907'
$ws.Range("E35").Value = '  -4.50%  '
$ws.Range("D36").Value = '62.02'
$ws.Range("E36").Value = '  +5.27%  '
$ws.Range("D37").Value = '6.06'
$ws.Range("E37").Value = '  +4.99%  '
$ws.Range("D38").Value = '40.80'
$ws.Range("E38").Value = '  -2.65%  '
$ws.Range("D40").Value = '0.997'
$ws.Range("E40").Value = '  -0.18%  '
$ws.Range("E41").Value = '  +0.11%  '
$ws.Range("E42").Value = '  +2.15%  '
$ws.Range("D43").Value = '2.95'
$ws.Range("E43").Value = '  +4.76%  '
$ws.Range("E44").Value = '  -1.11%  '
$ws.Range("D45").Value = '3.02'
$ws.Range("E45").Value = '  +2.24%  '
$ws.Range("E46").Value = '  +0.05%  '
$ws.Range("E47").Value = '  +6.67%  '
$ws.Range("B48").Value = 'BabyDogeCoin'
$ws.Range("C48").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D48").Value = '0.0₆0363'
$ws.Range("E48").Value = '  +10.84%  '
$ws.Range("B49").Value = 'Stacks'
$ws.Range("C49").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D49").Value = '3.07'
$ws.Range("E49").Value = '  +8.03%  '
$ws.Range("D50").Value = '3.39'
$ws.Range("E50").Value = '  -1.35%  '
$ws.Range("D51").Value = '2.08'
$ws.Range("E51").Value = '  -3.14%  '
